$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: the "Phone Number" column goes away, a new "SDT" column
#     appears in its place (shifted one column left) ---
# Old layout: D1=Salary, E1=Email, F1=Phone Number
# New layout: D1=Salary, E1=SDT,   F1=Email
$ws.Range("E1").Value = "SDT"
$ws.Range("F1").Value = "Email"

# --- Data rows: columns E (Email) and F (Phone) swap places ---
# Row 3 is the only one with data in BOTH columns (a genuine swap); the
# other rows only had an E value (Email) with F empty, so there the net
# effect is "value moves from E to F, E ends up blank". A temp-variable
# swap handles every case uniformly without clobbering anything.
function Swap-EF($row) {
    $eCell = $ws.Range("E$row")
    $fCell = $ws.Range("F$row")
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    # Writing a purely-numeric-looking string (like the phone number
    # "098271698 ") straight back through .Value would make Excel
    # re-interpret it as a number and drop the leading zero. Force it
    # to stay text with a leading apostrophe, then restore the cell's
    # original (unstyled) look so we don't leave a stray quote-prefix
    # style behind.
    if ($eVal -eq $null) { $eVal = "" }
    $fCell.Value = "'" + $eVal
    $fCell.Style = "Normal"

    if ($fVal -eq $null) {
        $eCell.Value = ""
    } else {
        $eCell.Value = "'" + $fVal
        $eCell.Style = "Normal"
    }
}

Swap-EF 3
Swap-EF 5
Swap-EF 24
Swap-EF 26
Swap-EF 28
Swap-EF 30
Swap-EF 33
Swap-EF 34

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = 31
$ws.Columns.Item(4).ColumnWidth = 11.7109375
$ws.Columns.Item(5).ColumnWidth = 13
$ws.Columns.Item(6).ColumnWidth = 25.7109375

# --- Sheet view: drop the frozen/scrolled topLeftCell, move selection to F1 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F1").Select()
